$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new sending-cluster category "M2" is inserted between FAPs and sCs.
# This pushes the former sCs block (rows 8-10) down to rows 11-13, and
# the new M2 block occupies rows 8-10. All data values are refreshed to
# match the updated (3 biological replicates) calculation.

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Vtn"
$ws.Cells.Item(2,3).Value = "Itgb8"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 12.75206033333332911
$ws.Cells.Item(2,8).Value = 38.25618099999999799
$ws.Cells.Item(2,9).Value = 0.15731223433819591
$ws.Cells.Item(2,10).Value = 0.15731223433819599
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.31166733333333341
$ws.Cells.Item(2,14).Value = 0.935002
$ws.Cells.Item(2,15).Value = 0.0414413620607491
$ws.Cells.Item(2,16).Value = 0.0414413620607491
$ws.Cells.Item(2,17).Value = 3.97440063859577819
$ws.Cells.Item(2,18).Value = 35.76960574736200016
$ws.Cells.Item(2,19).Value = 0.00651923325979458
$ws.Cells.Item(2,20).Value = 0.00651923325979459

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Vtn"
$ws.Cells.Item(3,3).Value = "Itgb8"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 12.75206033333332911
$ws.Cells.Item(3,8).Value = 38.25618099999999799
$ws.Cells.Item(3,9).Value = 0.15731223433819591
$ws.Cells.Item(3,10).Value = 0.15731223433819599
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.79458466666666716
$ws.Cells.Item(3,14).Value = 11.38375399999999971
$ws.Cells.Item(3,15).Value = 0.50455322140968761
$ws.Cells.Item(3,16).Value = 0.50455322140968761
$ws.Cells.Item(3,17).Value = 48.38877260927488777
$ws.Cells.Item(3,18).Value = 435.49895348347388335
$ws.Cells.Item(3,19).Value = 0.07937239460249244
$ws.Cells.Item(3,20).Value = 0.07937239460249246

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Vtn"
$ws.Cells.Item(4,3).Value = "Itgb8"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 12.75206033333332911
$ws.Cells.Item(4,8).Value = 38.25618099999999799
$ws.Cells.Item(4,9).Value = 0.15731223433819591
$ws.Cells.Item(4,10).Value = 0.15731223433819599
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 3.41443066666666706
$ws.Cells.Item(4,14).Value = 10.24329200000000029
$ws.Cells.Item(4,15).Value = 0.45400541652956322
$ws.Cells.Item(4,16).Value = 0.45400541652956328
$ws.Cells.Item(4,17).Value = 43.54102586531688246
$ws.Cells.Item(4,18).Value = 391.86923278785201319
$ws.Cells.Item(4,19).Value = 0.07142060647590891
$ws.Cells.Item(4,20).Value = 0.07142060647590892

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Vtn"
$ws.Cells.Item(5,3).Value = "Itgb8"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 20.35396833333334143
$ws.Cells.Item(5,8).Value = 61.06190500000001009
$ws.Cells.Item(5,9).Value = 0.2510910513649196
$ws.Cells.Item(5,10).Value = 0.2510910513649196
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.31166733333333341
$ws.Cells.Item(5,14).Value = 0.935002
$ws.Cells.Item(5,15).Value = 0.0414413620607491
$ws.Cells.Item(5,16).Value = 0.0414413620607491
$ws.Cells.Item(5,17).Value = 6.34366703320111291
$ws.Cells.Item(5,18).Value = 57.0930032988100109
$ws.Cells.Item(5,19).Value = 0.01040555516982778
$ws.Cells.Item(5,20).Value = 0.01040555516982778

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Vtn"
$ws.Cells.Item(6,3).Value = "Itgb8"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 20.35396833333334143
$ws.Cells.Item(6,8).Value = 61.06190500000001009
$ws.Cells.Item(6,9).Value = 0.2510910513649196
$ws.Cells.Item(6,10).Value = 0.2510910513649196
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.79458466666666716
$ws.Cells.Item(6,14).Value = 11.38375399999999971
$ws.Cells.Item(6,15).Value = 0.50455322140968761
$ws.Cells.Item(6,16).Value = 0.50455322140968761
$ws.Cells.Item(6,17).Value = 77.23485614348557249
$ws.Cells.Item(6,18).Value = 695.11370529137013818
$ws.Cells.Item(6,19).Value = 0.12668879883331549
$ws.Cells.Item(6,20).Value = 0.12668879883331549

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Vtn"
$ws.Cells.Item(7,3).Value = "Itgb8"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 20.35396833333334143
$ws.Cells.Item(7,8).Value = 61.06190500000001009
$ws.Cells.Item(7,9).Value = 0.2510910513649196
$ws.Cells.Item(7,10).Value = 0.2510910513649196
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 3.41443066666666706
$ws.Cells.Item(7,14).Value = 10.24329200000000029
$ws.Cells.Item(7,15).Value = 0.45400541652956322
$ws.Cells.Item(7,16).Value = 0.45400541652956328
$ws.Cells.Item(7,17).Value = 69.49721366569556835
$ws.Cells.Item(7,18).Value = 625.47492299126008675
$ws.Cells.Item(7,19).Value = 0.1139966973617763
$ws.Cells.Item(7,20).Value = 0.1139966973617763

# Row 8
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Vtn"
$ws.Cells.Item(8,3).Value = "Itgb8"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.33333333333333331
$ws.Cells.Item(8,7).Value = 0.004706
$ws.Cells.Item(8,8).Value = 0.014118
$ws.Cells.Item(8,9).Value = 0.00005805425597465
$ws.Cells.Item(8,10).Value = 0.00005805425597465
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.31166733333333341
$ws.Cells.Item(8,14).Value = 0.935002
$ws.Cells.Item(8,15).Value = 0.0414413620607491
$ws.Cells.Item(8,16).Value = 0.0414413620607491
$ws.Cells.Item(8,17).Value = 0.00146670647066667
$ws.Cells.Item(8,18).Value = 0.013200358236
$ws.Cells.Item(8,19).Value = 0.00000240584744101
$ws.Cells.Item(8,20).Value = 0.00000240584744101

# Row 9
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Vtn"
$ws.Cells.Item(9,3).Value = "Itgb8"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.33333333333333331
$ws.Cells.Item(9,7).Value = 0.004706
$ws.Cells.Item(9,8).Value = 0.014118
$ws.Cells.Item(9,9).Value = 0.00005805425597465
$ws.Cells.Item(9,10).Value = 0.00005805425597465
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.79458466666666716
$ws.Cells.Item(9,14).Value = 11.38375399999999971
$ws.Cells.Item(9,15).Value = 0.50455322140968761
$ws.Cells.Item(9,16).Value = 0.50455322140968761
$ws.Cells.Item(9,17).Value = 0.01785731544133333
$ws.Cells.Item(9,18).Value = 0.16071583897200001
$ws.Cells.Item(9,19).Value = 0.00002929146186855
$ws.Cells.Item(9,20).Value = 0.00002929146186855

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Vtn"
$ws.Cells.Item(10,3).Value = "Itgb8"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.33333333333333331
$ws.Cells.Item(10,7).Value = 0.004706
$ws.Cells.Item(10,8).Value = 0.014118
$ws.Cells.Item(10,9).Value = 0.00005805425597465
$ws.Cells.Item(10,10).Value = 0.00005805425597465
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 3.41443066666666706
$ws.Cells.Item(10,14).Value = 10.24329200000000029
$ws.Cells.Item(10,15).Value = 0.45400541652956322
$ws.Cells.Item(10,16).Value = 0.45400541652956328
$ws.Cells.Item(10,17).Value = 0.01606831071733333
$ws.Cells.Item(10,18).Value = 0.14461479645600001
$ws.Cells.Item(10,19).Value = 0.00002635694666509
$ws.Cells.Item(10,20).Value = 0.00002635694666509

# Row 11
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Vtn"
$ws.Cells.Item(11,3).Value = "Itgb8"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 47.95136699999999763
$ws.Cells.Item(11,8).Value = 143.85410100000001421
$ws.Cells.Item(11,9).Value = 0.59153866004090971
$ws.Cells.Item(11,10).Value = 0.59153866004090982
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.31166733333333341
$ws.Cells.Item(11,14).Value = 0.935002
$ws.Cells.Item(11,15).Value = 0.0414413620607491
$ws.Cells.Item(11,16).Value = 0.0414413620607491
$ws.Cells.Item(11,17).Value = 14.94487468257799989
$ws.Cells.Item(11,18).Value = 134.50387214320198837
$ws.Cells.Item(11,19).Value = 0.02451416778368571
$ws.Cells.Item(11,20).Value = 0.02451416778368572

# Row 12
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Vtn"
$ws.Cells.Item(12,3).Value = "Itgb8"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 47.95136699999999763
$ws.Cells.Item(12,8).Value = 143.85410100000001421
$ws.Cells.Item(12,9).Value = 0.59153866004090971
$ws.Cells.Item(12,10).Value = 0.59153866004090982
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 3.79458466666666716
$ws.Cells.Item(12,14).Value = 11.38375399999999971
$ws.Cells.Item(12,15).Value = 0.50455322140968761
$ws.Cells.Item(12,16).Value = 0.50455322140968761
$ws.Cells.Item(12,17).Value = 181.9555219639059942
$ws.Cells.Item(12,18).Value = 1637.59969767515394778
$ws.Cells.Item(12,19).Value = 0.29846273651201111
$ws.Cells.Item(12,20).Value = 0.29846273651201111

# Row 13
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Vtn"
$ws.Cells.Item(13,3).Value = "Itgb8"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 47.95136699999999763
$ws.Cells.Item(13,8).Value = 143.85410100000001421
$ws.Cells.Item(13,9).Value = 0.59153866004090971
$ws.Cells.Item(13,10).Value = 0.59153866004090982
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 3.41443066666666706
$ws.Cells.Item(13,14).Value = 10.24329200000000029
$ws.Cells.Item(13,15).Value = 0.45400541652956322
$ws.Cells.Item(13,16).Value = 0.45400541652956328
$ws.Cells.Item(13,17).Value = 163.72661799338800392
$ws.Cells.Item(13,18).Value = 1473.53956194049192163
$ws.Cells.Item(13,19).Value = 0.2685617557452129
$ws.Cells.Item(13,20).Value = 0.26856175574521302
